# Update Data by bot, scripted by HH
# Row 6 of Sheet1 refreshes the 605377 record from the 2019 annual report
# snapshot to the 2020 Q3 (three-quarter) report snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report date
$ws.Range("H6").Value = "2020-09-30 00:00:00"

# Numeric performance figures
$ws.Range("I6").Value = 1.13
$ws.Range("J6").Value = 1.12
$ws.Range("K6").Value = 1097936334.86
$ws.Range("L6").Value = 172812155.43
$ws.Range("M6").Value = 14.86
$ws.Range("N6").Value = -12.7063467786
$ws.Range("O6").Value = 55.13
$ws.Range("P6").Value = 8.168876543034999
$ws.Range("Q6").Value = -0.163291914846
$ws.Range("R6").Value = 25.3172256591

# These two figures are no longer populated for this report type
$ws.Range("S6").ClearContents()
$ws.Range("T6").ClearContents()

# Metadata describing which report/period this row represents.
# AB6 ("1") and AE6 ("2020") look numeric, so force text formatting first
# - otherwise COM auto-converts the assigned string into a Double, same as
# typing a bare number into a General-formatted cell in real Excel.
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value = "1"
$ws.Range("AB6").NumberFormat = "General"

$ws.Range("AC6").Value = "2020Q3"
$ws.Range("AD6").Value = "2020年 三季报"

$ws.Range("AE6").NumberFormat = "@"
$ws.Range("AE6").Value = "2020"
$ws.Range("AE6").NumberFormat = "General"

$ws.Range("AF6").Value = "三季报"
$ws.Range("AG6").Value = "2020-12-09 07:52:42"
